$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 now documents the (simplified) path to ShockMaster.xlsx, and the
# shock-logic subfolder "Electricity price logics" no longer exists, so the
# label in column A is renamed from "Electricity" to "ShockMaster".
$ws.Range("A10").Value = "ShockMaster"
$ws.Range("B10").Value = "C:\Users\loren\Documents\GitHub\SESAM\GT-IOA\Shocks\ShockMaster.xlsx"

# Reflect the saved cursor/selection position from the authored workbook.
$ws.Range("B17").Select()
